$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Créditos-aula: 4 -> 2 (numeric-looking text, force as text via quote prefix,
#     then restore the clean "text/wrap" style from the column so it matches the
#     sibling cells exactly instead of picking up a quote-prefixed style) ---
$ws.Range("B5").Value = "'2"
$ws.Range("B10").Copy()
$ws.Range("B5").PasteSpecial($xlPasteFormats)

$ws.Range("C5").Value = "'2"
$ws.Range("C10").Copy()
$ws.Range("C5").PasteSpecial($xlPasteFormats)

# --- Carga horária: 60 h -> 30 h (plain text, no special handling needed) ---
$ws.Range("B7").Value = "30 h"
$ws.Range("C7").Value = "30 h"

# --- Ativação: 01/01/2012 -> 01/01/2023 (date-looking text, force as text) ---
$ws.Range("B8").Value = "'01/01/2023"
$ws.Range("B10").Copy()
$ws.Range("B8").PasteSpecial($xlPasteFormats)

$ws.Range("C8").Value = "'01/01/2023"
$ws.Range("C10").Copy()
$ws.Range("C8").PasteSpecial($xlPasteFormats)

# --- Objetivos (docente responsável): Carlos Yujiro Shigue -> Emerson Gonçalves de Melo ---
$ws.Range("B10").Value = "7290967 - Emerson Gonçalves de Melo"
$ws.Range("C10").Value = "7290967 - Emerson Gonçalves de Melo"

# --- Objectives: new English objectives text (B11/C11 did not exist before,
#     copy the column's normal wrap-text format onto them first) ---
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial($xlPasteFormats)
$ws.Range("B11").Value = "Present the concepts of nanoscience and nanotechnology. The physical and chemical properties of materials on a nanometer scale are described by the laws of quantum mechanics, presenting in these dimensions different characteristics of materials on a macroscopic scale. The knowledge of this interdisciplinary area is fundamental in the formation of a researcher and/or a professional working in the area of materials."

$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial($xlPasteFormats)
$ws.Range("C11").Value = "Present the concepts of nanoscience and nanotechnology. The physical and chemical properties of materials on a nanometer scale are described by the laws of quantum mechanics, presenting in these dimensions different characteristics of materials on a macroscopic scale. The knowledge of this interdisciplinary area is fundamental in the formation of a researcher and/or a professional working in the area of materials."

# --- Programa resumido: Semestral -> 01/01/2023 (date-looking text) ---
$ws.Range("B13").Value = "'01/01/2023"
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial($xlPasteFormats)

$ws.Range("C13").Value = "'01/01/2023"
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial($xlPasteFormats)

# --- Short syllabus: (new B14/C14) ---
$ws.Range("B10").Copy()
$ws.Range("B14").PasteSpecial($xlPasteFormats)
$ws.Range("B14").Value = "Nanoscience and nanotechnology: principles and applications."

$ws.Range("C10").Copy()
$ws.Range("C14").PasteSpecial($xlPasteFormats)
$ws.Range("C14").Value = "Nanoscience and nanotechnology: principles and applications."

# --- Programa: (now shows the same docente name) ---
$ws.Range("B15").Value = "7290967 - Emerson Gonçalves de Melo"
$ws.Range("C15").Value = "7290967 - Emerson Gonçalves de Melo"

# --- Syllabus: (new B16/C16, full English syllabus) ---
$ws.Range("B10").Copy()
$ws.Range("B16").PasteSpecial($xlPasteFormats)
$ws.Range("B16").Value = "Conceptualization: nanoscience and nanotechnology. Low-dimensional systems. Quantum Confinement. Chemical bonds: molecules and clusters. Electronic and structural properties. Synthesis and fabrication of materials at the nanometer scale: bottom-up and top-down techniques. Fullerenes and carbon nanotubes. Molecular self-organization and supramolecular systems. Quantum wires and dots. Magnetic nanoparticles. Characterization techniques: X-ray diffraction, scattering and absorption, scanning tunneling microscopy (STM), atomic force microscopy (AFM), transmission electron microscopy. Transport properties: ballistic transport, quantum conductance, Coulomb blocking. Molecular devices. Diffusive transport. Nanomagnetism: magnetic order, superparamagnetism and spintronics. Applications."

$ws.Range("C10").Copy()
$ws.Range("C16").PasteSpecial($xlPasteFormats)
$ws.Range("C16").Value = "Conceptualization: nanoscience and nanotechnology. Low-dimensional systems. Quantum Confinement. Chemical bonds: molecules and clusters. Electronic and structural properties. Synthesis and fabrication of materials at the nanometer scale: bottom-up and top-down techniques. Fullerenes and carbon nanotubes. Molecular self-organization and supramolecular systems. Quantum wires and dots. Magnetic nanoparticles. Characterization techniques: X-ray diffraction, scattering and absorption, scanning tunneling microscopy (STM), atomic force microscopy (AFM), transmission electron microscopy. Transport properties: ballistic transport, quantum conductance, Coulomb blocking. Molecular devices. Diffusive transport. Nanomagnetism: magnetic order, superparamagnetism and spintronics. Applications."

# --- Método: (new second docente responsável) ---
$ws.Range("B18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
